$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them (e.g. "11.60" -> 11.6,
# losing the trailing zero / exact text representation required by the sheet).

$ws.Range('D2').Value = '34.862.39'
$ws.Range('E2').Value = '  -0.56%  '

$ws.Range('D3').Value = '1.840.17'
$ws.Range('E3').Value = '  +1.71%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.69'
$ws.Range('E5').Value = '  -0.39%  '

$ws.Range('E6').Value = '  +1.03%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '39.99'
$ws.Range('E8').Value = '  -0.70%  '

$ws.Range('E9').Value = '  +0.68%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0686'
$ws.Range('E10').Value = '  +0.27%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0980'
$ws.Range('E11').Value = '  -1.84%  '

$ws.Range('D12').Value = '2.109.90'

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.60'
$ws.Range('E13').Value = '  +4.89%  '

$ws.Range('D14').Value = '1.847.48'
$ws.Range('E14').Value = '  +3.14%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.673'
$ws.Range('E15').Value = '  +1.36%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.64'
$ws.Range('E16').Value = '  -0.36%  '

$ws.Range('D17').Value = '34.866.25'
$ws.Range('E17').Value = '  -0.41%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.82'
$ws.Range('E18').Value = '  +0.07%  '

$ws.Range('D19').Value = '0.0₃0786'
$ws.Range('E19').Value = '  -0.45%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '239.92'
$ws.Range('E20').Value = '  +0.92%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.15'
$ws.Range('E21').Value = '  +1.72%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.68'
$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('E24').Value = '  +0.94%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '170.94'
$ws.Range('E25').Value = '  -0.70%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.79'
$ws.Range('E26').Value = '  -1.03%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.45'
$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('E28').Value = '  +2.04%  '

$ws.Range('E29').Value = '  -3.42%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0551'
$ws.Range('E31').Value = '  -1.03%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.94'
$ws.Range('E32').Value = '  -4.55%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.96'
$ws.Range('E33').Value = '  -1.62%  '

$ws.Range('E34').Value = '  +8.33%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.22'
$ws.Range('E35').Value = '  +6.69%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.46'
$ws.Range('E36').Value = '  +13.66%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.695'
$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('E38').Value = '  +7.63%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '90.23'
$ws.Range('E39').Value = '  -1.84%  '

$ws.Range('D40').Value = '1.348.12'
$ws.Range('E40').Value = '  +2.65%  '

$ws.Range('E41').Value = '  +0.14%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '14.79'
$ws.Range('E42').Value = '  +2.29%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.29'
$ws.Range('E43').Value = '  +1.26%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').Value = '  -2.49%  '

$ws.Range('E45').Value = '  +0.17%  '

$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0523'
$ws.Range('E46').Value = '  +2.29%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.25'
$ws.Range('E47').Value = '  -0.81%  '

$ws.Range('D48').Value = '2.023.45'
$ws.Range('E48').Value = '  +1.88%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.47'
$ws.Range('E49').Value = '  +23.09%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.06%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0668'
$ws.Range('E51').Value = '  +0.54%  '
